$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 590, shifting existing rows (590 onward) down by 2
$ws.Rows("590:591").Insert()

# Row 590: new data record (Primera)
$ws.Range("A590").Value() = 6
$ws.Range("B590").Value() = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C590").Value() = "Metropolitana"
$ws.Range("D590").Value() = "2022-03-24"
$ws.Range("E590").Value() = 13
$ws.Range("F590").Value() = 100112023
$ws.Range("G590").Value() = "Brócoli"
$ws.Range("H590").Value() = "Sin especificar"
$ws.Range("I590").Value() = "Primera"
$ws.Range("J590").Value() = 11700
$ws.Range("K590").Value() = 700
$ws.Range("L590").Value() = 900
$ws.Range("M590").Value() = 804
$ws.Range("N590").Value() = "$/unidad"
$ws.Range("O590").Value() = "Región Metropolitana"
$ws.Range("P590").Value() = 804
$ws.Range("Q590").Value() = 1
$ws.Range("R590").Value() = "Hortaliza"

# Row 591: new data record (Segunda)
$ws.Range("A591").Value() = 6
$ws.Range("B591").Value() = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C591").Value() = "Metropolitana"
$ws.Range("D591").Value() = "2022-03-24"
$ws.Range("E591").Value() = 13
$ws.Range("F591").Value() = 100112023
$ws.Range("G591").Value() = "Brócoli"
$ws.Range("H591").Value() = "Sin especificar"
$ws.Range("I591").Value() = "Segunda"
$ws.Range("J591").Value() = 6300
$ws.Range("K591").Value() = 500
$ws.Range("L591").Value() = 600
$ws.Range("M591").Value() = 546
$ws.Range("N591").Value() = "$/unidad"
$ws.Range("O591").Value() = "Región Metropolitana"
$ws.Range("P591").Value() = 546
$ws.Range("Q591").Value() = 1
$ws.Range("R591").Value() = "Hortaliza"
